# Update the crypto price/volume table (rows 2-51) with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "51.706.44"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +1.39%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.031.46"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "380.49"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.51%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "102.86"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.65%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "36.84"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  +1.30%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "3.514.57"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  -0.74%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.035.80"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("E17").Value = "  -3.58%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "10.53"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -15.73%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "51.701.20"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.43%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "3.07"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.94%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.50"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +0.86%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0963"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E23").Value = "  +0.93%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "268.49"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  -2.62%  "
$ws.Range("E26").Value = "  +2.43%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.63"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +8.66%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.173"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("E29").Value = "  -0.06%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "26.27"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("E31").Value = "  +0.97%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "10.29"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.67%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "34.14"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.81%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "50.53"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -0.08%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0448"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +6.45%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.294"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +13.63%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "17.06"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +2.99%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.86"
$cell.ClearFormats()
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("E44").Value = "  +5.83%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "123.80"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +4.46%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "21.85"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("E47").Value = "  +3.82%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.41"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +4.22%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.036.10"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +1.40%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.333.90"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("E51").Value = "  +0.27%  "
